$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update recalculated "Adj Close" (F) values and a few other cells in existing rows ---
$ws.Range("F2").Value2 = 99.01917299999999
$ws.Range("F4").Value2 = 100.795403
$ws.Range("F5").Value2 = 102.649696
$ws.Range("F6").Value2 = 103.19622
$ws.Range("F7").Value2 = 101.488319
$ws.Range("F8").Value2 = 103.137657
$ws.Range("F10").Value2 = 104.181923
$ws.Range("F11").Value2 = 103.508522
$ws.Range("F15").Value2 = 103.157188
$ws.Range("F16").Value2 = 104.640617
$ws.Range("F17").Value2 = 104.894371
$ws.Range("F18").Value2 = 104.406395
$ws.Range("F19").Value2 = 104.142883
$ws.Range("F20").Value2 = 103.762283
$ws.Range("F21").Value2 = 104.59182
$ws.Range("F23").Value2 = 105.986115
$ws.Range("F25").Value2 = 103.628044
$ws.Range("F26").Value2 = 102.796364
$ws.Range("F28").Value2 = 102.307144
$ws.Range("F29").Value2 = 101.162346
$ws.Range("F31").Value2 = 101.778763
$ws.Range("F32").Value2 = 100.829681
$ws.Range("F33").Value2 = 99.401123
$ws.Range("F34").Value2 = 100.174103
$ws.Range("F35").Value2 = 98.22699
$ws.Range("F38").Value2 = 98.794487
$ws.Range("F39").Value2 = 99.068443
$ws.Range("F40").Value2 = 99.518547
$ws.Range("F43").Value2 = 99.940315
$ws.Range("F45").Value2 = 99.773567
$ws.Range("F46").Value2 = 99.871651
$ws.Range("F47").Value2 = 100.116867
$ws.Range("F48").Value2 = 103.569511
$ws.Range("F50").Value2 = 102.098213
$ws.Range("F51").Value2 = 104.069763
$ws.Range("F52").Value2 = 103.255638
$ws.Range("F55").Value2 = 102.990791
$ws.Range("F56").Value2 = 104.364014
$ws.Range("F57").Value2 = 104.364014
$ws.Range("F59").Value2 = 102.343422
$ws.Range("F60").Value2 = 102.529793
$ws.Range("F62").Value2 = 102.794617
$ws.Range("F63").Value2 = 104.334595
$ws.Range("F64").Value2 = 104.825287
$ws.Range("F65").Value2 = 105.346466
$ws.Range("F66").Value2 = 106.447815
$ws.Range("F67").Value2 = 106.723145
$ws.Range("F68").Value2 = 105.002281
$ws.Range("F69").Value2 = 105.21862
$ws.Range("F71").Value2 = 104.284447
$ws.Range("F77").Value2 = 102.661903
$ws.Range("F79").Value2 = 105.179276
$ws.Range("F80").Value2 = 104.058266
$ws.Range("F81").Value2 = 103.025749
$ws.Range("F82").Value2 = 104.687607
$ws.Range("F83").Value2 = 101.668991
$ws.Range("F84").Value2 = 104.202576
$ws.Range("F85").Value2 = 104.784233
$ws.Range("F86").Value2 = 103.7491
$ws.Range("F89").Value2 = 101.590134
$ws.Range("F90").Value2 = 102.575958
$ws.Range("F91").Value2 = 103.660378
$ws.Range("F92").Value2 = 102.792839
$ws.Range("F93").Value2 = 101.728142
$ws.Range("F94").Value2 = 101.422531
$ws.Range("F96").Value2 = 100.377556
$ws.Range("F97").Value2 = 99.66773999999999
$ws.Range("F98").Value2 = 99.312859
$ws.Range("F99").Value2 = 99.598732
$ws.Range("F100").Value2 = 99.10581999999999
$ws.Range("F101").Value2 = 98.85936700000001
$ws.Range("F104").Value2 = 101.530968
$ws.Range("F108").Value2 = 101.217644
$ws.Range("F109").Value2 = 99.71517900000001
$ws.Range("F110").Value2 = 100.881554
$ws.Range("F111").Value2 = 100.743172
$ws.Range("F112").Value2 = 101.039719
$ws.Range("F113").Value2 = 100.041374
$ws.Range("F116").Value2 = 101.415321
$ws.Range("F117").Value2 = 102.117119
$ws.Range("F120").Value2 = 102.136894
$ws.Range("F121").Value2 = 102.245636
$ws.Range("F122").Value2 = 101.978745
$ws.Range("F123").Value2 = 102.413673
$ws.Range("F124").Value2 = 100.565254
$ws.Range("F129").Value2 = 98.20117999999999
$ws.Range("F130").Value2 = 98.330017
$ws.Range("F132").Value2 = 99.935654
$ws.Range("F138").Value2 = 100.797935
$ws.Range("F139").Value2 = 100.827667
$ws.Range("F141").Value2 = 100.272629
$ws.Range("F145").Value2 = 99.162575
$ws.Range("F151").Value2 = 96.096504
$ws.Range("F152").Value2 = 96.59343699999999
$ws.Range("F153").Value2 = 95.00324999999999
$ws.Range("F160").Value2 = 91.95210299999999
$ws.Range("F162").Value2 = 94.953568
$ws.Range("F163").Value2 = 94.327438
$ws.Range("F164").Value2 = 94.635536
$ws.Range("F165").Value2 = 94.734909
$ws.Range("G201").Value2 = 63724600
$ws.Range("D202").Value2 = 82.739998
$ws.Range("E202").Value2 = 82.769997
$ws.Range("F202").Value2 = 82.769997
$ws.Range("G202").Value2 = 87696900

# --- Append two new rows of historical data (2023-10-20 and 2023-10-23) ---
$rA = $ws.Range("A203")
$rA.NumberFormat = "@"
$rA.Value2 = "2023-10-20"
$rA.ClearFormats()
$ws.Range("B203").Value2 = 82.989998
$ws.Range("C203").Value2 = 83.540001
$ws.Range("D203").Value2 = 82.769997
$ws.Range("E203").Value2 = 83.239998
$ws.Range("F203").Value2 = 83.239998
$ws.Range("G203").Value2 = 52162600

$rA = $ws.Range("A204")
$rA.NumberFormat = "@"
$rA.Value2 = "2023-10-23"
$rA.ClearFormats()
$ws.Range("B204").Value2 = 82.989998
$ws.Range("C204").Value2 = 82.894997
$ws.Range("D204").Value2 = 82.58000199999999
$ws.Range("E204").Value2 = 82.58429700000001
$ws.Range("F204").Value2 = 82.58429700000001
$ws.Range("G204").Value2 = 3530145

